# Automated map update:
# - Row 41 (Caso -347, AVALOS /ALT/ 1610) is removed entirely; all following
#   rows shift up by one.
# - A brand-new case (Caso 6472, DUMONT SANTOS 3744) is appended as the new
#   last row of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 41; Excel shifts rows 42..89 up to 41..88 automatically.
$ws.Rows.Item(41).Delete()

# Append the new record as row 89.
$newRow = 89

# Columns that hold text values in this table (everything except
# Attachments/Coordenada_X/Coordenada_Y), force them to be stored as text so
# numeric-looking values (case numbers, comunas, OT numbers, dates) keep
# their original text representation instead of being auto-converted.
$textCols = @("A","B","C","D","E","F","G","H","J","K","L","O","P")
foreach ($col in $textCols) {
    $ws.Range("$col$newRow").NumberFormat = "@"
}

$ws.Range("A$newRow").Value = "6472"
$ws.Range("B$newRow").Value = "7/24/2025"
$ws.Range("C$newRow").Value = "DUMONT, SANTOS 3744"
$ws.Range("D$newRow").Value = "15"
$ws.Range("E$newRow").Value = "808509381"
$ws.Range("F$newRow").Value = "AYKO"
$ws.Range("G$newRow").Value = "Pendiente"
$ws.Range("H$newRow").Value = "Columna inclinada"
$ws.Range("I$newRow").Value = 1
$ws.Range("J$newRow").Value = "Aplomo"
$ws.Range("K$newRow").Value = "Sin equipos"
$ws.Range("L$newRow").Value = "Poste"
$ws.Range("M$newRow").Value = -58.448576
$ws.Range("N$newRow").Value = -34.585794
$ws.Range("O$newRow").Value = "Colegiales"
$ws.Range("P$newRow").Value = "Capital Norte"
